# Auto-generated edit script for Recommandations + Top_YTD sheets
# Applies the BRVM automatic update (GitHub Actions) diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# ---- Sheet "Recommandations" (A2:G48) ----
$data1 = @(
    @("CFAO MOTORS CI", 0, 4, 3695, 960, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI", 0, 4, 3565, 805, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 8, 3311.34, 104, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 4, 2800, 700, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2615, 680, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2496.58, 615.64, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2385, 590, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1591.17, 395.92, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1389.08, 348.8, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1318.23, 329.93, "🟡 Observer", "➖ Neutre"),
    @("Marché", 0, 1, 805, 805, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 551.23, 138.17, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 533.64, 134.13, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 505.8, 128.07, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 497.08, 125.86, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 495.13, 122.66, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 432.23, 107.77, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 2, 419.67, 210.16, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 377.89, 94.42, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 1, 251.95, 251.95, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE             (**)", 0, 1, 251.86, 251.86, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                (**)", 0, 1, 191.84, 191.84, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 1, 191.6, 191.6, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 3, 1, 13.8, -6.22, "🟢 Achat", "✅ Renforcer"),
    @("ORAGROUP TOGO (ORGT)", 2, 0, 12.5, 6.93, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 1, 0, 6.98, 6.98, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 4.26, 4.26, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA BN (BOAB)", 1, 0, 4.01, 4.01, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 1, 0, 3.44, 3.44, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 1, 0, 3.34, 3.34, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 1, 0, 3.23, 3.23, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 1, 2, 2.5, -2.11, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA ML (BOAM)", 1, 0, 2.23, 2.23, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 1, 1, 0.37, 6.25, "🟡 Observer", "👀 À surveiller"),
    @("SMB CI (SMBC)", 1, 1, 0.19, 7.07, "🟡 Observer", "👀 À surveiller"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 1, 1, 0.09, 3.14, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA NG (BOAN)", 1, 1, -0.9399999999999999, 4, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 0, 1, -1.13, -1.13, "🟡 Observer", "➖ Neutre"),
    @("VIVO ENERGY CI (SHEC)", 0, 1, -1.3, -1.3, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 0, 1, -1.96, -1.96, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING SN (TTLS)", 0, 1, -2.34, -2.34, "🟡 Observer", "➖ Neutre"),
    @("SOGB CI (SOGC)", 0, 1, -2.47, -2.47, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 1, 2, -5.77, -6.04, "🟡 Observer", "👀 À surveiller"),
    @("SICOR CI (SICC)", 0, 1, -7.38, -7.38, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 1, 2, -9.35, -7.47, "🟡 Observer", "👀 À surveiller"),
    @("UNILEVER CI (UNLC)", 0, 3, -21.13, -6.67, "🔴 Vente", "⚠️ Risque de décrochage")
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = $data1[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws1.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}

# ---- Sheet "Top_YTD" (A2:B11) ----
$data2 = @(
    @("BRVM - SERVICES PUBLICS", 7951751.82),
    @("CFAO MOTORS CI", 1093250.25),
    @("UNIWAX CI", 960224.23),
    @("SETAO CI", 408869),
    @("NEI-CEDA CI", 322366.92),
    @("BRVM - AUTRES SECTEURS", 274837.75),
    @("AIR LIQUIDE CI", 234879.5),
    @("BRVM - DISTRIBUTION", 61292.15),
    @("BRVM - TRANSPORT", 39917.36),
    @("BRVM - AGRICULTURE", 33947.48)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = $data2[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws2.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}

Write-Host "Update complete"
